# netCrypto.xlsx — "Add files via upload" re-upload edit.
#
# The underlying change recorded in the commit is a refreshed snapshot of
# the same workbook: the author's local folder name changed (BACKUP - Copy
# -> 01062025), Excel stamped a new co-authoring revision GUID and window
# position on its own when the file was re-saved, and the one substantive
# data edit is cell T2 on SheetName1, which was bumped from 10 to 60593.
#
# Window chrome / absolute-path / revision-GUID bookkeeping are Excel-
# internal metadata that Excel itself rewrites on every save and that the
# object model does not expose knobs for (there's no Range/Workbook
# property that maps onto <x15ac:absPath>, <xr:revisionPtr>, or the
# <workbookView xWindow/yWindow> attributes) - so this script focuses on
# the one cell that actually carries information: T2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# Substantive data change from the diff: T2 10 -> 60593
$ws.Range("T2").Value = 60593
